$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H38").Value = 142.72728
$ws.Range("I38").Value = 142.72728
$ws.Range("K38").Value = 428.18184
$ws.Range("M38").Value = -56.18184000000002

$ws.Range("H46").Value = 8008.5
$ws.Range("I46").Value = 8008.5
$ws.Range("K46").Value = 24025.5
$ws.Range("M46").Value = -23906.5

$ws.Range("H51").Value = 11149.8
$ws.Range("I51").Value = 9999
$ws.Range("K51").Value = 9999
$ws.Range("M51").Value = -9515

$ws.Range("H60").Value = 8008.5
$ws.Range("I60").Value = 8008.5
$ws.Range("K60").Value = 24025.5
$ws.Range("M60").Value = -23541.5

$ws.Range("H70").Value = 1992.3334
$ws.Range("I70").Value = 1860.75
$ws.Range("J70").Value = 2142.7144
$ws.Range("K70").Value = 5582.25
$ws.Range("L70").Value = 6428.1432
$ws.Range("M70").Value = -5312.25
$ws.Range("N70").Value = -6968.1432

$ws.Range("H73").Value = 1992.3334
$ws.Range("I73").Value = 1860.75
$ws.Range("J73").Value = 2142.7144
$ws.Range("K73").Value = 5582.25
$ws.Range("L73").Value = 6428.1432
$ws.Range("M73").Value = -4646.25
$ws.Range("N73").Value = -8300.143199999999

$ws.Range("H80").Value = 949.8182
$ws.Range("I80").Value = 596.8570999999999
$ws.Range("J80").Value = 1567.5
$ws.Range("K80").Value = 1790.5713
$ws.Range("L80").Value = 4702.5
$ws.Range("M80").Value = -792.5712999999998
$ws.Range("N80").Value = -6698.5

$ws.Range("H83").Value = 949.8182
$ws.Range("I83").Value = 596.8570999999999
$ws.Range("J83").Value = 1567.5
$ws.Range("K83").Value = 5371.7139
$ws.Range("L83").Value = 14107.5
$ws.Range("M83").Value = -379.7138999999997
$ws.Range("N83").Value = -24091.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3867.5745
$ws.Range("I32").Value = 1214.0975
$ws.Range("K32").Value = 1214.0975
$ws.Range("M32").Value = -927.0975000000001

$ws.Range("H61").Value = 2034.6666
$ws.Range("I61").Value = 2100.8333
$ws.Range("J61").Value = 1770
$ws.Range("K61").Value = 2100.8333
$ws.Range("L61").Value = 1770
$ws.Range("M61").Value = -1888.8333
$ws.Range("N61").Value = -2194

$ws.Range("H62").Value = 55113
$ws.Range("J62").Value = 50000
$ws.Range("L62").Value = 50000
$ws.Range("N62").Value = -51248

$ws.Range("H65").Value = 55113
$ws.Range("J65").Value = 50000
$ws.Range("L65").Value = 150000
$ws.Range("N65").Value = -156240

$ws.Range("H76").Value = 28858.334
$ws.Range("J76").Value = 28858.334
$ws.Range("L76").Value = 28858.334
$ws.Range("N76").Value = -29534.334

$ws.Range("H79").Value = 28858.334
$ws.Range("J79").Value = 28858.334
$ws.Range("L79").Value = 28858.334
$ws.Range("N79").Value = -31198.334

$ws.Range("H136").Value = 2034.6666
$ws.Range("I136").Value = 2100.8333
$ws.Range("J136").Value = 1770
$ws.Range("K136").Value = 6302.499899999999
$ws.Range("L136").Value = 5310
$ws.Range("M136").Value = -3752.499899999999
$ws.Range("N136").Value = -10410

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H9").Value = 0
$ws.Range("J9").Value = 0
$ws.Range("L9").Value = 0
$ws.Range("N9").Value = $null

$ws.Range("H76").Value = 19989.166
$ws.Range("J76").Value = 19989.166
$ws.Range("L76").Value = 19989.166
$ws.Range("N76").Value = -20619.166

$ws.Range("H79").Value = 19989.166
$ws.Range("J79").Value = 19989.166
$ws.Range("L79").Value = 19989.166
$ws.Range("N79").Value = -22173.166

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1689.6
$ws.Range("I31").Value = 1324.5
$ws.Range("J31").Value = 1933
$ws.Range("K31").Value = 1324.5
$ws.Range("L31").Value = 1933
$ws.Range("M31").Value = -1029.5
$ws.Range("N31").Value = -2523

$ws.Range("H34").Value = 1689.6
$ws.Range("I34").Value = 1324.5
$ws.Range("J34").Value = 1933
$ws.Range("K34").Value = 1324.5
$ws.Range("L34").Value = 1933
$ws.Range("M34").Value = -1122.5
$ws.Range("N34").Value = -2337

$ws.Range("H58").Value = 1717.375
$ws.Range("I58").Value = 1760.2667
$ws.Range("K58").Value = 1760.2667
$ws.Range("M58").Value = -1557.2667

$ws.Range("H88").Value = 38244
$ws.Range("J88").Value = 38244
$ws.Range("L88").Value = 38244
$ws.Range("N88").Value = -39056

$ws.Range("H91").Value = 38244
$ws.Range("J91").Value = 38244
$ws.Range("L91").Value = 38244
$ws.Range("N91").Value = -41052

$ws.Range("H136").Value = 1717.375
$ws.Range("I136").Value = 1760.2667
$ws.Range("K136").Value = 5280.800099999999
$ws.Range("M136").Value = -2730.800099999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 766.55554
$ws.Range("I2").Value = 566.6667
$ws.Range("J2").Value = 866.5
$ws.Range("K2").Value = 3400.0002
$ws.Range("L2").Value = 5199
$ws.Range("M2").Value = -3287.0002
$ws.Range("N2").Value = -5425

$ws.Range("H4").Value = 1603.2084
$ws.Range("I4").Value = 1459.8334
$ws.Range("J4").Value = 2033.3334
$ws.Range("K4").Value = 4379.5002
$ws.Range("L4").Value = 6100.0002
$ws.Range("M4").Value = -4267.5002
$ws.Range("N4").Value = -6324.0002

$ws.Range("H17").Value = 220.33333
$ws.Range("J17").Value = 220.33333
$ws.Range("L17").Value = 660.99999
$ws.Range("N17").Value = -998.99999

$ws.Range("H108").Value = 569.3333
$ws.Range("I108").Value = 569.3333
$ws.Range("K108").Value = 1707.9999
$ws.Range("M108").Value = 1172.0001

$ws.Range("H113").Value = 674.2
$ws.Range("I113").Value = 342.25
$ws.Range("J113").Value = 1053.5714
$ws.Range("K113").Value = 1026.75
$ws.Range("L113").Value = 3160.7142
$ws.Range("M113").Value = 1143.25
$ws.Range("N113").Value = -7500.7142

$ws.Range("H129").Value = 1728
$ws.Range("I129").Value = 1125.6
$ws.Range("J129").Value = 2732
$ws.Range("K129").Value = 3376.8
$ws.Range("L129").Value = 8196
$ws.Range("M129").Value = 1623.2
$ws.Range("N129").Value = -18196

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H11").Value = 978904.7
$ws.Range("I11").Value = 553868
$ws.Range("J11").Value = 2084000
$ws.Range("K11").Value = 553868
$ws.Range("L11").Value = 2084000
$ws.Range("M11").Value = -553729
$ws.Range("N11").Value = -2084278

$ws.Range("H12").Value = 5401.5
$ws.Range("I12").Value = 5401.5
$ws.Range("J12").Value = 0
$ws.Range("K12").Value = 5401.5
$ws.Range("L12").Value = 0
$ws.Range("M12").Value = -5261.5
$ws.Range("N12").Value = $null

$ws.Range("H14").Value = 13584
$ws.Range("J14").Value = 15251
$ws.Range("L14").Value = 15251
$ws.Range("N14").Value = -15587

$ws.Range("H80").Value = 1846.4
$ws.Range("I80").Value = 1846.4
$ws.Range("K80").Value = 1846.4
$ws.Range("M80").Value = -848.4000000000001

$ws.Range("H83").Value = 1846.4
$ws.Range("I83").Value = 1846.4
$ws.Range("K83").Value = 9232
$ws.Range("M83").Value = -4240

$ws.Range("H102").Value = 12496.25
$ws.Range("I102").Value = 3621.818
$ws.Range("K102").Value = 3621.818
$ws.Range("M102").Value = -1999.818

$ws.Range("H122").Value = 3307.9
$ws.Range("J122").Value = 3425.8572
$ws.Range("L122").Value = 10277.5716
$ws.Range("N122").Value = -15177.5716

$ws.Range("H132").Value = 2486.2856
$ws.Range("I132").Value = 2501.6
$ws.Range("J132").Value = 2448
$ws.Range("K132").Value = 7504.799999999999
$ws.Range("L132").Value = 7344
$ws.Range("M132").Value = -4974.799999999999
$ws.Range("N132").Value = -12404

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H17").Value = 16699
$ws.Range("I17").Value = 10502.667
$ws.Range("J17").Value = 25993.5
$ws.Range("K17").Value = 10502.667
$ws.Range("L17").Value = 25993.5
$ws.Range("M17").Value = -10332.667
$ws.Range("N17").Value = -26333.5

$ws.Range("H18").Value = 10000
$ws.Range("J18").Value = 10000
$ws.Range("L18").Value = 10000
$ws.Range("N18").Value = -10344

$ws.Range("H40").Value = 0
$ws.Range("I40").Value = 0
$ws.Range("K40").Value = 0
$ws.Range("M40").Value = $null

$ws.Range("H64").Value = 24615.666
$ws.Range("J64").Value = 24615.666
$ws.Range("L64").Value = 24615.666
$ws.Range("N64").Value = -25065.666

$ws.Range("H67").Value = 24615.666
$ws.Range("J67").Value = 24615.666
$ws.Range("L67").Value = 24615.666
$ws.Range("N67").Value = -26175.666

$ws.Range("H93").Value = 1500
$ws.Range("I93").Value = 1500
$ws.Range("K93").Value = 1500
$ws.Range("M93").Value = -252

$ws.Range("H132").Value = 2530
$ws.Range("I132").Value = 1990
$ws.Range("J132").Value = 2800
$ws.Range("K132").Value = 5970
$ws.Range("L132").Value = 8400
$ws.Range("M132").Value = -3440
$ws.Range("N132").Value = -13460

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 5427.143
$ws.Range("I62").Value = 5427.143
$ws.Range("K62").Value = 5427.143
$ws.Range("M62").Value = -4803.143

$ws.Range("H65").Value = 5427.143
$ws.Range("I65").Value = 5427.143
$ws.Range("K65").Value = 27135.715
$ws.Range("M65").Value = -24015.715

$ws.Range("H126").Value = 44531.727
$ws.Range("I126").Value = 38985
$ws.Range("K126").Value = 116955
$ws.Range("M126").Value = -114485

$ws.Range("H132").Value = 4645.5
$ws.Range("I132").Value = 0
$ws.Range("J132").Value = 4645.5
$ws.Range("K132").Value = 0
$ws.Range("L132").Value = 13936.5
$ws.Range("N132").Value = -18996.5
$ws.Range("M132").Value = $null
